$wb = $excel.ActiveWorkbook

# Rename the compound "L-Serine" -> "serine".
# The compound (col J) and compoundId (col K) columns hold this value for
# rows 2-9 on the original/cor_pct/cor_abs sheets, and cell B2 on the
# "total" sheet. Re-stamping the font after the value write nudges the
# cell back onto the workbook's implicit default style (matching the
# lightweight style shift seen alongside the rename).
$sheetNames = @("original", "cor_pct", "cor_abs")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($row = 2; $row -le 9; $row++) {
        foreach ($col in @("J", "K")) {
            $cell = $ws.Range($col + $row)
            $cell.Value = "serine"
            $cell.Font.Name = "Calibri"
            $cell.Font.Size = 11
        }
    }
}

$totalWs = $wb.Worksheets.Item("total")
$totalCell = $totalWs.Range("B2")
$totalCell.Value = "serine"
$totalCell.Font.Name = "Calibri"
$totalCell.Font.Size = 11

# Move the active tab/selection to the "total" sheet, landing on B2 - the
# cell that was just edited - mirroring where the user ended up.
[void]$totalWs.Activate()
[void]$totalWs.Range("B2").Select()

# Restore a sensible selection on the other sheets touched above (they
# settle on the compound/compoundId pair that was just renamed).
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $b2 = $ws.Range("B2")
    $k2 = $ws.Range("K2")
    $union = $excel.Union($b2, $k2)
    [void]$union.Select()
    [void]$k2.Activate()
}

[void]$totalWs.Activate()
